$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Bemessungsleistung PL (column D, rows 2-10) from 10000000 to 1000000.
# Column G holds formulas (=15*D{row}) and will recalc automatically.
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 4).Value = 1000000
}

# Update the active selection to match the saved view (D17 instead of D16).
$ws.Range("D17").Select()
